$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# "Summary: " -> "Summary:" (teal 21918C) + " " (purple 440154)
$teal = 9212193    # RGB(0x21, 0x91, 0x8C) packed as BGR for PowerPoint's RGB()
$purple = 5505348  # RGB(0x44, 0x01, 0x54) packed as BGR for PowerPoint's RGB()

$summaryLabel = $tr.Characters(1, 8)
$summaryLabel.Font.Color.RGB = $teal

$summarySpace = $tr.Characters(9, 1)
$summarySpace.Font.Color.RGB = $purple

# Color the "What's New?" and "Class Connection:" run labels teal as well
$full = $tr.Text

$whatIdx = $full.IndexOf("What")
$whatRange = $tr.Characters($whatIdx + 1, 11)
$whatRange.Font.Color.RGB = $teal

$classIdx = $full.IndexOf("Class Connection")
$classRange = $tr.Characters($classIdx + 1, 17)
$classRange.Font.Color.RGB = $teal
